$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "testing MSE" (B) and "validation MSE" (C) values for rows 2-20.
$data = @(
    @(0.01652577235857634, 0.03347919488030672),
    @(0.01927675928009841, 0.03636165096783885),
    @(0.02159953269702994, 0.04077332017719157),
    @(0.02285051840486269, 0.04349266310446182),
    @(0.02520421430104667, 0.04831450552228939),
    @(0.02294172383415028, 0.04624239128749408),
    @(0.02421330825042563, 0.0483049256956185),
    @(0.02588327534669511, 0.04839857970646459),
    @(0.02660737802725797, 0.05635550058836643),
    @(0.02819156066213258, 0.0575588030636747),
    @(0.02833601731714437, 0.05688414398712061),
    @(0.02656732681116657, 0.0507410754640036),
    @(0.02849029058414042, 0.05428956764375461),
    @(0.02768074674960684, 0.05029263491849386),
    @(0.02803702131581991, 0.05413640313146839),
    @(0.02877949356325445, 0.05222792004919209),
    @(0.02832971697418312, 0.05245502903561287),
    @(0.02746171320502883, 0.05117971212170454),
    @(0.02807359794915116, 0.05398631178995315)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
}
